$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the letter / correct-answer pairs that changed between the two
# condition files (condition 1 n1-back data).
$ws.Range("B5").Value = "H"
$ws.Range("C5").Value = 0

$ws.Range("B6").Value = "H"
$ws.Range("C6").Value = 1

$ws.Range("B11").Value = "K"
$ws.Range("C11").Value = 1

$ws.Range("C12").Value = 0

$ws.Range("B35").Value = "T"
$ws.Range("C35").Value = 0

$ws.Range("B40").Value = "A"
$ws.Range("C40").Value = 1

# Move the active selection to match the author's saved view.
$ws.Range("D6").Select() | Out-Null
